$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert a new column for "C.Reynolds" right before "T.Williams" (col J / 10),
    # shifting every later player column one to the right.
    $ws.Columns.Item(10).Insert()
    $ws.Cells.Item(1, 10).Value = "C.Reynolds"
    $ws.Cells.Item(2, 10).Value = "n"

    # Append a new trailing column for "S.Zylstra" right after "C.Taumoepeau"
    # (now col X / 24), inserting one column past the last used column so it
    # picks up the header formatting from its left neighbor.
    $ws.Columns.Item(25).Insert()
    $ws.Cells.Item(1, 25).Value = "S.Zylstra"
    $ws.Cells.Item(2, 25).Value = "n"
}
